$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: replace the two "Palace Square / Naberezhnye Chelny-ish" pair with the London pair ---
$ws.Range("A3").Value = "212 Baker St, Marylebone, London"
$ws.Range("B3").Value = "51.5242146, -0.1584717"

# --- Row 4: new short strings, B4 loses its wrap-text styling ---
$ws.Range("A4").Value = "jhbkvcz"
$ws.Range("B4").Value = "ыпчваспрмглш"
$ws.Range("B4").WrapText = $false

# --- Row 5: A5 becomes a single space, B5 is removed entirely ---
$ws.Range("A5").Value = " "
$ws.Range("A5").WrapText = $false
$ws.Range("B5").Clear()

# --- Row 6: brand-new numeric row ---
$ws.Range("A6").Value = 1234132456549
$ws.Range("A6").NumberFormat = "#,##0"
$ws.Range("B6").Value = 0.8765432
$ws.Range("B6").WrapText = $false

# --- Row 7: A7 becomes "%", B7 is removed entirely ---
$ws.Range("A7").Value = "%"
$ws.Range("A7").WrapText = $false
$ws.Range("B7").Clear()

# --- column widths (closest achievable values given pixel-quantized COM width units) ---
$ws.Range("A1").EntireColumn.ColumnWidth = 35.83333333333333
$ws.Range("B1").EntireColumn.ColumnWidth = 33.33333333333333

# the old file also defined an explicit (default-width) column group for C:AMK;
# drop it so only columns A and B carry an explicit width, matching the edit.
# (EntireColumn.Delete keeps the used range intact, unlike ClearFormats.)
$tailCols = $ws.Range($ws.Cells.Item(1, 3), $ws.Cells.Item(1, 1025))
$tailCols.EntireColumn.Delete()

# --- selection moves to B6 ---
$ws.Range("B6").Select()
